$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "San Diego Comic-Con 2016 (PS16)"
$ws.Range("A2").Value = "Chandra, Flamecaller"
$ws.Range("A3").Value = "Gideon, Ally of Zendikar"
$ws.Range("A4").Value = "Jace, Unraveler of Secrets"
$ws.Range("A5").Value = "Liliana, the Last Hope"
$ws.Range("A6").Value = "Nissa, Voice of Zendikar"
